$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update price (D) and 1h-volume-change (E) columns with the latest scrape.
# D-column prices are stored as TEXT (some contain multiple "." thousand
# separators, e.g. "52.055.51", which are not valid numbers) -- force the
# cell to text before assigning so Excel does not auto-convert numeric-looking
# strings (e.g. "342.77") into real numbers, then drop back to the default
# "Normal" style so no stray number-format style lingers on the cell.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "52.055.51"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +4.63%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.783.44"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +5.00%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "342.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.70%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "115.57"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.26%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.550"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.88%  "

$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("E9").Value = "  +4.94%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.15"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.84%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0858"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.95%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "20.00"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.13%  "

$ws.Range("E13").Value = "  +1.49%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.64"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.69%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.218.43"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.10%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.782.20"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.02%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "51.948.87"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.53%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.879"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.28%  "

$ws.Range("E19").Value = "  +9.90%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.02"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.79%  "

$ws.Range("E21").Value = "  -1.01%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0982"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.27%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "277.05"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.85%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.14"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.38%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.78"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +7.99%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.67"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.75%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.01%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.20"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.25%  "

$ws.Range("E29").Value = "  +0.82%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.140"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.93%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.79"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.52%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "50.24"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.24%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.71"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.25%  "

$ws.Range("E34").Value = "  -0.23%  "

$ws.Range("E35").Value = "  +0.03%  "

$ws.Range("E36").Value = "  +3.20%  "

$ws.Range("E37").Value = "  -1.26%  "

$ws.Range("E38").Value = "  +0.60%  "

$ws.Range("E39").Value = "  +2.75%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0383"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +10.90%  "

$ws.Range("E41").Value = "  +26.55%  "

$ws.Range("E42").Value = "  +3.68%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.35"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.04%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "23.41"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.70%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "127.12"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.41%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.072.79"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.28%  "

$ws.Range("E47").Value = "  -0.34%  "

$ws.Range("E48").Value = "  +0.49%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.57"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.47%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.897"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +15.83%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.84"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.23%  "
